# Fix line-break placeholders (double spaces) in the "Extreme Low Flow Method"
# labels for rows D and F on the ExtremeFlows sheet, and leave the selection
# on the next empty row (A12) as it was when the author last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExtremeFlows")

$ws.Range("A7").Value = "F. Lowest consecutive flows in  Reclamation's ensembles and traces"
$ws.Range("A5").Value = "D. 85%, 65%, and 50% of  2000 to 2018 average flow"

$ws.Range("A12").Select()
